$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 19:03"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1218968
$ws.Range("C4").Value = 6133
$ws.Range("E4").Value = 959431
$ws.Range("G4").Value = 838
$ws.Range("H4").Value = 70759

# Row 61 - Luxemburgo
$ws.Range("B61").Value = 3840
$ws.Range("C61").Value = 12
$ws.Range("D61").Value = 3412
$ws.Range("E61").Value = 332
$ws.Range("F61").Value = 22

# Row 71 - Irak
$ws.Range("B71").Value = 2431
$ws.Range("C71").Value = 85
$ws.Range("D71").Value = 1571
$ws.Range("E71").Value = 758
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 102

# Row 91 - Republica de Yibuti
$ws.Range("B91").Value = 1120
$ws.Range("C91").Value = 4
$ws.Range("D91").Value = 745
$ws.Range("E91").Value = 373
